$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - copy H1's formatting (bold,
# centered, bordered header style) so the new cells share the same style
# index as the rest of the header row, then overwrite the text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I and J, rows 2-4
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 6
